# Developer Guide: Fix spelling error
# ------------------------------------
# 1) Correct the misspelled "CrearCommand" shape label to "ClearCommand"
#    on the Design slide.
# 2) The deck was also re-saved by PowerPoint, which re-stamped the
#    "datetimeFigureOut" date field cached on the slide master and every
#    slide layout (8/7/2018 -> 9/21/18).

$p = $ppt.ActivePresentation

# --- 1) Spelling fix: CrearCommand -> ClearCommand ------------------------
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame) {
            if ($shape.TextFrame.TextRange.Text -eq "CrearCommand") {
                $shape.TextFrame.TextRange.Text = "ClearCommand"
            }
        }
    }
}

# --- 2) Refresh the cached date field text --------------------------------
$oldDate = "8/7/2018"
$newDate = "9/21/18"

function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        if ($shape.Name -like "Date Placeholder*" -and $shape.HasTextFrame) {
            if ($shape.TextFrame.TextRange.Text -eq $oldDate) {
                $shape.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholders $layouts.Item($li).Shapes
}

Write-Output "edit complete"
